# Digikey/Regulators.xlsx: rename the two worksheets to match the file
# name ("Reg-SMD" -> "Regulators-SMD", "Reg-TH" -> "Regulators-TH") and
# make the SMD sheet the active/selected tab instead of the TH sheet.

$wb = $excel.ActiveWorkbook

$wsSMD = $wb.Worksheets.Item("Reg-SMD")
$wsSMD.Name = "Regulators-SMD"

$wsTH = $wb.Worksheets.Item("Reg-TH")
$wsTH.Name = "Regulators-TH"

# Switch the selected/active tab from Reg-TH to Reg-SMD (first sheet).
$wsSMD.Activate()
